# Daily attendance processing - 2025-12-06 16:54:30
# Normalizes the ordering of names in the "Recorded By" column (G) so that
# "System" (and "admin@admin.com") is always listed before the human
# reviewer's e-mail address / lower-case "system" token.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match replacements to apply to column G ("Recorded By")
$replacements = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "system, System, backup@backdoor.com" = "System, system, backup@backdoor.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is column index 7
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
